$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Values are written with a leading apostrophe to force Excel to keep
# them as literal text (matching the original inlineStr cell type)
# instead of auto-converting numeric-looking strings to numbers;
# the Style reset back to "Normal" clears the quote-prefix marker
# afterwards so no stray cell formatting is introduced.

$ws.Range('D2').Value = "'64.680.93"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -3.03%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.171.51"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  -1.99%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.14%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'598.66"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -1.26%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'151.51"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -4.45%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.14%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'3.172.62"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  -1.90%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'  -3.75%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = "'  -5.63%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'5.55"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -2.72%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.474"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -6.63%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'0.0000259"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -6.14%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'36.85"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -5.79%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'3.694.56"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -1.82%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'64.701.77"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -3.04%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'3.185.58"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -1.36%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.35%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'7.01"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -5.36%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'479.13"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -6.43%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'14.78"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -3.06%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'0.713"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  -3.21%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'7.72"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -4.33%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'13.80"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -5.91%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'84.00"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -1.29%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'0.999"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.28%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'2.91"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -2.93%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  -6.18%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'2.25"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -6.03%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  +0.23%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  +11.61%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('B32').Value = "'FirstDigitalUSD"
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = "'1.00"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.05%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('B33').Value = "'Stacks"
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = "'2.72"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -8.51%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'26.82"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -5.21%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  -6.53%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'6.11"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -6.38%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'54.62"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -2.08%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'3.20"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +4.04%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  -5.47%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'458.75"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -10.21%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = "'Kaspa"
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = "'0.125"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -3.17%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = "'VeChain"
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = "'0.0402"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -5.14%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'8.47"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -3.25%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'2.41"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -2.38%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'2.868.09"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -0.37%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'0.274"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -8.59%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'27.03"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -5.40%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  +0.05%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('B49').Value = "'Stellar"
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = "'0.116"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.82%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = "'ThetaToken"
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'2.33"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -4.02%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'119.83"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -2.26%  "
$ws.Range('E51').Style = 'Normal'
